# Add a new worksheet "policy taxonomy" at the end of the workbook
# (after the current last sheet, "Sheet1") containing a small table that
# lists password-policy classes, their requirements and example
# passwords.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$ws.Name = "policy taxonomy"

# Populate the cells in the same order the original author filled them in
# (this also controls the order new entries land in the shared-string table).
$ws.Range("A2").Value = "comp8"
$ws.Range("A5").Value = "2word16"
$ws.Range("A4").Value = "3class12"
$ws.Range("A3").Value = "basic8 (1class8)"
$ws.Range("C1").Value = "Example passwords"
$ws.Range("A1").Value = "Policy"
$ws.Range("B5").Value = "be at least 16 characters long and include at least two letter sequences that are separated by a non-letter sequence."
$ws.Range("B1").Value = "The password needs to ..."
$ws.Range("B4").Value = "be at least 12 characters long and include three different character classes (upper, lower, digits, symbol)"
$ws.Range("B3").Value = "be at least 8 characters long"
$ws.Range("B2").Value = "bet at least 8 characters long, include at least one character from each character class, and not include a dictionary word"
$ws.Range("C3").Value = "password -- monkey123 -- qwerasdf"
$ws.Range("C2").Value = "P@ssw0rd -- !M0nkey1 -- LGtjj{Rd;w1u\"
$ws.Range("C4").Value = "Password1234 -- 2MonkeysBite -- NfJidl2kdils"
$ws.Range("C5").Value = "password.unlocks -- 1-Monkey-Bites -- qwer.asdf.zxcvb.1234"

# Header row: bold
$ws.Range("A1:C1").Font.Bold = $true

# First column body cells: explicit (non-bold) font, matches the new style
# entry added to styles.xml
$ws.Range("A2:A5").Font.Name = "Calibri"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 14.43
$ws.Columns.Item(2).ColumnWidth = 99.45
$ws.Columns.Item(3).ColumnWidth = 41.17

# Select the whole table and make this the active sheet/tab
$ws.Range("A1:C5").Select() | Out-Null
$ws.Activate()
